# Update the workbook to match the new scrape:
#  - refresh the "Förändrad" (Changed) date column C for all existing
#    data rows (2..127) from 45204 to 45205
#  - append five new logging-notification rows (128..132) for HÖRBY

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. bump the "Förändrad" column for every existing data row -----------
$ws.Range("C2:C127").Value = 45205

# --- 2. append the new rows -------------------------------------------------
$newRows = @(
    @{ Row = 128; A = "A 47609-2023"; B = 45203; C = 45205; Area = 2.4 },
    @{ Row = 129; A = "A 47647-2023"; B = 45203; C = 45205; Area = 1.1 },
    @{ Row = 130; A = "A 47612-2023"; B = 45203; C = 45205; Area = 2.4 },
    @{ Row = 131; A = "A 47637-2023"; B = 45203; C = 45205; Area = 2 },
    @{ Row = 132; A = "A 47641-2023"; B = 45203; C = 45205; Area = 4.8 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = "SKÅNE LÄN"
    $ws.Cells.Item($row, 5).Value = "HÖRBY"

    # G = Area (ha)
    $ws.Cells.Item($row, 7).Value = $r.Area

    # H..Q = species-count columns, all zero for these new entries
    $ws.Range($ws.Cells.Item($row, 8), $ws.Cells.Item($row, 17)).Value = 0

    # B & C keep the date number format used throughout the sheet
    $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 3)).NumberFormat = "YYYY-MM-DD"

    # R = Artnamn, wrap-text styled but left empty like the sibling rows
    $ws.Cells.Item($row, 18).WrapText = $true
}

# rows 128-131 carry an explicit 15pt row height (matching the rest of the
# sheet); row 132, being the new last row, is left at the implicit default
# the same way row 127 previously was.
$ws.Range("A127:A131").RowHeight = 15
